# Apply updated odds values to Sheet1, matching the target OOXML diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("G2").Value = 1.65
$ws.Range("I2").Value = 5.5
$ws.Range("J2").Value = 2.25
$ws.Range("U2").Value = 1.91
$ws.Range("V2").Value = 1.91

# Row 3
$ws.Range("Q3").Value = 2.6
$ws.Range("R3").Value = 1.48

# Row 4
$ws.Range("G4").Value = 1.65
$ws.Range("H4").Value = 3.4
$ws.Range("I4").Value = 6.25
$ws.Range("M4").Value = 1.08
$ws.Range("N4").Value = 8
$ws.Range("U4").Value = 2.25
$ws.Range("V4").Value = 1.57
$ws.Range("X4").Value = 6.5
$ws.Range("AB4").Value = 41
$ws.Range("AJ4").Value = 21
$ws.Range("AN4").Value = 3.4
$ws.Range("AP4").Value = 26
$ws.Range("AQ4").Value = 29
$ws.Range("AR4").Value = 67
$ws.Range("AS4").Value = 251
$ws.Range("AU4").Value = 10
$ws.Range("AZ4").Value = 151
$ws.Range("BA4").Value = 201

# Row 13
$ws.Range("O13").Value = 1.18
$ws.Range("P13").Value = 4.5
$ws.Range("Q13").Value = 1.62
$ws.Range("R13").Value = 2.25
